# "Ajuste busca de e-mail, criação de arquivo .xlsx com o resultado das buscas"
#
# On the "Base" sheet, column A held a mix of numeric voucher codes and
# text codes. The search/lookup logic was tightened so that every row in
# A1:A7 is now formatted as Text (so codes like "1234" can never be
# silently read back as a number), and the two rows that used to hold the
# raw numeric vouchers (1234 / 12345) are corrected to the proper textual
# voucher code "ASNUAS" (matching what the other rows already used).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Base")

# Format column A (rows 1-7) as Text. This introduces the new cell style
# (numFmtId 49 == "@") applied to A1:A7.
$ws.Range("A1:A7").NumberFormat = "@"

# Fix the two rows that still held raw numeric vouchers; they should hold
# the "ASNUAS" text code like row 4 already does.
$ws.Range("A2").Value = "ASNUAS"
$ws.Range("A3").Value = "ASNUAS"
